$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) cells are plain text (not numbers), e.g. "26.366.71" or
# "306.02" -- force text format on just the cells we are rewriting so Excel
# keeps them as literal strings instead of coercing to numeric values.
# (Each contiguous block is set separately -- a single comma-joined multi-area
# Range only applies NumberFormat to its first area.)
$ws.Range("D2:D18").NumberFormat = "@"
$ws.Range("D20:D26").NumberFormat = "@"
$ws.Range("D28:D34").NumberFormat = "@"
$ws.Range("D36:D46").NumberFormat = "@"
$ws.Range("D48:D51").NumberFormat = "@"

# --- Row 45/46 swap: Quant <-> EnergySwap (coin name, link, price, volume) ---
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("B46").Value = 'Quant'
$ws.Range("C46").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'

# --- Price (column D) updates ---
$ws.Range("D2").Value = '26.366.71'
$ws.Range("D3").Value = '1.774.30'
$ws.Range("D4").Value = '1.003'
$ws.Range("D5").Value = '1.002'
$ws.Range("D6").Value = '306.02'
$ws.Range("D7").Value = '0.4230'
$ws.Range("D8").Value = '0.3605'
$ws.Range("D9").Value = '0.07143'
$ws.Range("D10").Value = '0.8361'
$ws.Range("D11").Value = '20.35'
$ws.Range("D12").Value = '1.784.38'
$ws.Range("D13").Value = '6.444'
$ws.Range("D14").Value = '5.229'
$ws.Range("D15").Value = '0.06839'
$ws.Range("D16").Value = '1.003'
$ws.Range("D17").Value = '78.66'
$ws.Range("D18").Value = '0.000008663'
$ws.Range("D20").Value = '14.88'
$ws.Range("D21").Value = '26.388.39'
$ws.Range("D22").Value = '5.063'
$ws.Range("D23").Value = '11.00'
$ws.Range("D24").Value = '2.016.72'
$ws.Range("D25").Value = '152.21'
$ws.Range("D26").Value = '1.795'
$ws.Range("D28").Value = '5.051'
$ws.Range("D29").Value = '113.95'
$ws.Range("D30").Value = '1.812'
$ws.Range("D31").Value = '0.08820'
$ws.Range("D32").Value = '0.7235'
$ws.Range("D33").Value = '1.116'
$ws.Range("D34").Value = '4.310'
$ws.Range("D36").Value = '2.736'
$ws.Range("D37").Value = '1.092'
$ws.Range("D38").Value = '0.05125'
$ws.Range("D39").Value = '0.01879'
$ws.Range("D40").Value = '0.1608'
$ws.Range("D41").Value = '0.4896'
$ws.Range("D42").Value = '2.616'
$ws.Range("D43").Value = '6.346'
$ws.Range("D44").Value = '7.958'
$ws.Range("D45").Value = '10.22'
$ws.Range("D46").Value = '104.49'
$ws.Range("D48").Value = '1.637'
$ws.Range("D49").Value = '0.06178'
$ws.Range("D50").Value = '0.4451'
$ws.Range("D51").Value = '1.713'

# --- Volume(1h) (column E) updates ---
$ws.Range("E2").Value = '  -3.12%  '
$ws.Range("E3").Value = '  -1.98%  '
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("E5").Value = '  +0.07%  '
$ws.Range("E7").Value = '  +0.48%  '
$ws.Range("E8").Value = '  +1.20%  '
$ws.Range("E9").Value = '  +0.03%  '
$ws.Range("E10").Value = '  -1.72%  '
$ws.Range("E11").Value = '  +0.84%  '
$ws.Range("E12").Value = '  -2.79%  '
$ws.Range("E13").Value = '  +1.15%  '
$ws.Range("E14").Value = '  -1.64%  '
$ws.Range("E15").Value = '  -0.14%  '
$ws.Range("E16").Value = '  +0.18%  '
$ws.Range("E17").Value = '  -2.75%  '
$ws.Range("E18").Value = '  -1.00%  '
$ws.Range("E19").Value = '  +0.02%  '
$ws.Range("E20").Value = '  -1.61%  '
$ws.Range("E21").Value = '  -3.38%  '
$ws.Range("E22").Value = '  -0.89%  '
$ws.Range("E23").Value = '  +1.21%  '
$ws.Range("E24").Value = '  -1.76%  '
$ws.Range("E25").Value = '  -0.92%  '
$ws.Range("E26").Value = '  -8.91%  '
$ws.Range("E27").Value = '  -0.52%  '
$ws.Range("E28").Value = '  -0.57%  '
$ws.Range("E29").Value = '  +0.31%  '
$ws.Range("E30").Value = '  +7.58%  '
$ws.Range("E31").Value = '  -1.15%  '
$ws.Range("E32").Value = '  -1.97%  '
$ws.Range("E33").Value = '  +0.85%  '
$ws.Range("E34").Value = '  -2.81%  '
$ws.Range("E35").Value = '  +0.03%  '
$ws.Range("E36").Value = '  -6.93%  '
$ws.Range("E37").Value = '  +1.66%  '
$ws.Range("E38").Value = '  -1.02%  '
$ws.Range("E39").Value = '  -1.45%  '
$ws.Range("E40").Value = '  -1.52%  '
$ws.Range("E41").Value = '  -1.57%  '
$ws.Range("E42").Value = '  -3.61%  '
$ws.Range("E43").Value = '  +1.08%  '
$ws.Range("E44").Value = '  -2.70%  '
$ws.Range("E45").Value = '  -0.12%  '
$ws.Range("E46").Value = '  -0.56%  '
$ws.Range("E47").Value = '  +0.06%  '
$ws.Range("E48").Value = '  +2.61%  '
$ws.Range("E49").Value = '  -2.96%  '
$ws.Range("E50").Value = '  -2.75%  '
$ws.Range("E51").Value = '  +1.81%  '
